$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "EIA Table 2.13" -> "EIA Table 2.14" (formula references update
#    automatically when the sheet is renamed).
# ---------------------------------------------------------------------------
$eia = $wb.Worksheets.Item("EIA Table 2.13")
$eia.Name = "EIA Table 2.14"

$about = $wb.Worksheets.Item("About")
$calc  = $wb.Worksheets.Item("Calculations")
$tcamrb = $wb.Worksheets.Item("TCAMRB")

# ---------------------------------------------------------------------------
# 2. EIA Table 2.14 sheet: drop the 2009 row, append a 2020 row.
# ---------------------------------------------------------------------------
$eia.Rows.Item(5).Delete()

# Append the new 2020 data row (after delete, row 15 holds the notes row;
# insert a fresh row 15 above it for the 2020 figures).
$eia.Rows.Item(15).Insert()
$eia.Range("A15").Value = 2020
$eia.Range("B15").Value = 57001240
$eia.Range("C15").Value = 9855106
$eia.Range("D15").Value = 4447623
$eia.Range("E15").Value = 4279573
$eia.Range("F15").Value = 61448863
$eia.Range("G15").Value = 14134679

# Match styling of the row above (2019, now row 14) for the new 2020 row.
$eia.Range("A14:G14").Copy()
$eia.Range("A15:G15").PasteSpecial(-4122) # xlPasteFormats

$eia.Range("A2").Value2 = "https://www.eia.gov/electricity/annual/"
$eia.Range("A2").Value2 = "and Mexico, 2010-2020 (Megawatthours)"

$eia.Range("A1").Value2 = "Table 2.14.  Electric Power Industry - U.S. Electricity Imports from and Electricity Exports to Canada"

# ---------------------------------------------------------------------------
# 3. Calculations sheet: update the summed ranges to reflect the shifted
#    EIA Table 2.14 data (2010-2012 instead of 2011-2013).
# ---------------------------------------------------------------------------
$calc.Range("B21").Formula = "=SUM('EIA Table 2.14'!D5:E7)"
$calc.Range("B22").Formula = "=SUM('EIA Table 2.14'!B5:C7)"

$calc.Range("B23").Select()

# ---------------------------------------------------------------------------
# 4. About sheet updates.
# ---------------------------------------------------------------------------
$about.Range("B12").Value = 2021
$about.Range("B13").Value2 = "Electric Power Annual 2021 (with data for 2020)"
$about.Range("B15").Value2 = "Table 2.14"

$about.Hyperlinks.Add($about.Range("B14"), "https://www.eia.gov/electricity/annual/") | Out-Null

$wb.Save()
